# Refresh the "cryptos" price/volume snapshot (Price column D, Volume(1h) column E)
# on Sheet1, rows 2-51, with the latest scraped figures.
#
# The Price column stores values as plain text (e.g. "27.056.29", "1.006",
# "0.07350") rather than numbers, so every Price cell is written with a
# leading single-quote to force Excel to keep it as text - this avoids the
# default Value-setter auto-converting a string like "1.006" into the
# number 1.006 (which would also silently drop significant trailing zeros,
# e.g. "0.07350" -> 0.0735).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''27.056.29'
$ws.Range('E2').Value = '  -0.23%  '
$ws.Range('D3').Value = '''1.829.46'
$ws.Range('E3').Value = '  +0.19%  '
$ws.Range('E4').Value = '  -0.71%  '
$ws.Range('D5').Value = '''312.24'
$ws.Range('E5').Value = '  -0.14%  '
$ws.Range('D6').Value = '''1.006'
$ws.Range('E6').Value = '  -0.36%  '
$ws.Range('D7').Value = '''0.4604'
$ws.Range('E7').Value = '  -0.73%  '
$ws.Range('D8').Value = '''0.3698'
$ws.Range('E8').Value = '  +1.61%  '
$ws.Range('D9').Value = '''0.07350'
$ws.Range('E9').Value = '  +0.49%  '
$ws.Range('D10').Value = '''0.8716'
$ws.Range('E10').Value = '  -0.38%  '
$ws.Range('D11').Value = '''0.07925'
$ws.Range('E11').Value = '  +3.81%  '
$ws.Range('D12').Value = '''19.80'
$ws.Range('E12').Value = '  -1.63%  '
$ws.Range('D13').Value = '''1.780.48'
$ws.Range('E13').Value = '  -4.02%  '
$ws.Range('D14').Value = '''5.347'
$ws.Range('E14').Value = '  -0.17%  '
$ws.Range('D15').Value = '''6.562'
$ws.Range('E15').Value = '  +1.57%  '
$ws.Range('D16').Value = '''91.74'
$ws.Range('E16').Value = '  -1.07%  '
$ws.Range('E17').Value = '  -0.09%  '
$ws.Range('D18').Value = '''0.000008874'
$ws.Range('E18').Value = '  +2.38%  '
$ws.Range('D19').Value = '''1.006'
$ws.Range('E19').Value = '  -0.35%  '
$ws.Range('D20').Value = '''14.71'
$ws.Range('E20').Value = '  +1.51%  '
$ws.Range('D21').Value = '''27.030.88'
$ws.Range('E21').Value = '  -1.96%  '
$ws.Range('D22').Value = '''5.123'
$ws.Range('E22').Value = '  -1.92%  '
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').Value = '''1.948.51'
$ws.Range('E24').Value = '  -6.66%  '
$ws.Range('D25').Value = '''152.41'
$ws.Range('E25').Value = '  +0.77%  '
$ws.Range('D26').Value = '''1.846'
$ws.Range('E26').Value = '  -1.84%  '
$ws.Range('D27').Value = '''18.49'
$ws.Range('E27').Value = '  +0.33%  '
$ws.Range('D28').Value = '''2.068'
$ws.Range('E28').Value = '  -0.37%  '
$ws.Range('D29').Value = '''5.090'
$ws.Range('E29').Value = '  -0.82%  '
$ws.Range('D30').Value = '''115.18'
$ws.Range('E30').Value = '  -0.94%  '
$ws.Range('D31').Value = '''0.08868'
$ws.Range('E31').Value = '  -0.41%  '
$ws.Range('D32').Value = '''2.975'
$ws.Range('E32').Value = '  +0.52%  '
$ws.Range('D33').Value = '''0.7335'
$ws.Range('E33').Value = '  -0.10%  '
$ws.Range('D34').Value = '''4.435'
$ws.Range('E34').Value = '  -1.00%  '
$ws.Range('E35').Value = '  -1.81%  '
$ws.Range('D36').Value = '''2.457'
$ws.Range('E36').Value = '  -2.88%  '
$ws.Range('E37').Value = '  -1.86%  '
$ws.Range('D38').Value = '''0.05240'
$ws.Range('E38').Value = '  -0.62%  '
$ws.Range('D39').Value = '''0.01938'
$ws.Range('E39').Value = '  +0.25%  '
$ws.Range('D40').Value = '''2.945'
$ws.Range('E40').Value = '  +0.28%  '
$ws.Range('D41').Value = '''7.150'
$ws.Range('E41').Value = '  -1.51%  '
$ws.Range('D42').Value = '''0.5150'
$ws.Range('E42').Value = '  -1.23%  '
$ws.Range('D43').Value = '''0.1631'
$ws.Range('E43').Value = '  -0.05%  '
$ws.Range('D44').Value = '''0.8587'
$ws.Range('E44').Value = '  -15.06%  '
$ws.Range('D45').Value = '''8.232'
$ws.Range('E45').Value = '  -0.79%  '
$ws.Range('D46').Value = '''0.4825'
$ws.Range('E46').Value = '  -0.70%  '
$ws.Range('D47').Value = '''1.007'
$ws.Range('E47').Value = '  -0.33%  '
$ws.Range('D48').Value = '''10.18'
$ws.Range('E48').Value = '  -2.17%  '
$ws.Range('D49').Value = '''102.33'
$ws.Range('E49').Value = '  -1.47%  '
$ws.Range('E50').Value = '  -1.03%  '
$ws.Range('D51').Value = '''0.06223'
$ws.Range('E51').Value = '  -0.92%  '
